# Daily attendance processing - 2026-01-12 11:35:14
# Swap "System, <email>" -> "<email>, System" in column G (Recorded By)
# of the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val -like "System, *") {
        $rest = $val.Substring(8)
        if ($rest -notlike "*,*") {
            $cell.Value2 = "$rest, System"
        }
    }
}
